$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete columns AA:AB (design summer temperature, design winter temperature)
# This shifts the old AC column ("duration (minutes)") left to become the new AA.
$ws.Range("AA1:AB1").EntireColumn.Delete()

# Delete rows 4:5 (the PECO and BGE preset rows)
$ws.Range("A4:A5").EntireRow.Delete()

# Update the selected cell to match the target view state
$ws.Range("D3").Select()
